# Append a new "TAGS:" section (bold heading + four plain tag lines)
# after the existing "Role 2:" paragraph at the end of the document body.

$d = $word.ActiveDocument

# --- "TAGS:" bold heading paragraph -----------------------------------
# Inserting right after the last paragraph ("Role 2:") naturally inherits
# that paragraph's bold run formatting, which is exactly what we want for
# the "TAGS:" heading.
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$tagsPara = $d.Paragraphs.Last
$tagsPara.Range.Text = "TAGS:"

# --- Plain (non-bold) tag lines ----------------------------------------
$tagLines = @(
    "1 Slaaviq -- Russian Christmas",
    "1 Agayumayaraq, Agayumaciq -- Christianity",
    "1 Kass'alugpiaq -- Russian Orthodox",
    "1 Agayulirtet -- Priests"
)

foreach ($line in $tagLines) {
    $prevPara = $d.Paragraphs.Last
    $prevPara.Range.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Last
    # Select + ClearFormatting strips the inherited bold (from the
    # preceding paragraph) so the new line is plain text, matching the
    # source document's un-bolded tag paragraphs.
    $newPara.Range.Select()
    $word.Selection.ClearFormatting()
    $word.Selection.Text = $line
}

Write-Host "Appended TAGS section with" $tagLines.Count "tag lines."
